$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold, centered, bordered - same style as H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-23.
$data = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 5, 5),
    @(5, 4, 5),
    @(6, 6, 7),
    @(7, 12, 13),
    @(8, 5, 8),
    @(9, 8, 8),
    @(10, 6, 8),
    @(11, 8, 8),
    @(12, 7, 7),
    @(13, 7, 8),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 5, 6),
    @(17, 7, 8),
    @(18, 7, 7),
    @(19, 1, 3),
    @(20, 1, 2),
    @(21, 1, 1),
    @(22, 4, 5),
    @(23, 8, 8)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $i0 = $entry[1]
    $if = $entry[2]
    $ws.Cells.Item($row, 9).Value = $i0
    $ws.Cells.Item($row, 10).Value = $if
}

Write-Host "I0 and IF columns added"
